$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: dates (as Excel serial numbers) and hours
$data = @(
    @{ Row = 8;  Date = 41554; Hours = 4 },
    @{ Row = 9;  Date = 41555; Hours = 2 },
    @{ Row = 10; Date = 41556; Hours = 6 },
    @{ Row = 11; Date = 41557; Hours = 1 }
)

foreach ($item in $data) {
    $r = $item.Row
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $item.Date
    # Reuse the same number format as the existing date cells (e.g. A7)
    $ws.Cells.Item(7, 1).Copy() | Out-Null
    $cellA.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $item.Hours
}
$excel.CutCopyMode = 0

# Update the selected cell to match the post-edit state
$ws.Range("A12").Select()
